$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.542.90'
$ws.Range("E2").Value = '  -1.28%  '
$ws.Range("D3").Value = '1.577.72'
$ws.Range("E3").Value = '  -3.36%  '
$ws.Range("E4").Value = '  +0.27%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '205.87'
$ws.Range("E5").Value = '  -2.60%  '
$ws.Range("E6").Value = '  -3.42%  '
$ws.Range("E7").Value = '  +0.30%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.06'
$ws.Range("E8").Value = '  -5.86%  '
$ws.Range("E9").Value = '  -2.34%  '
$ws.Range("E10").Value = '  -3.89%  '
$ws.Range("E11").Value = '  -2.22%  '
$ws.Range("D12").Value = '1.800.55'
$ws.Range("E12").Value = '  -3.41%  '
$ws.Range("D13").Value = '1.548.75'
$ws.Range("E13").Value = '  -5.21%  '
$ws.Range("E14").Value = '  -4.86%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.525'
$ws.Range("E15").Value = '  -7.09%  '
$ws.Range("D16").Value = '27.506.67'
$ws.Range("E16").Value = '  -1.45%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.61'
$ws.Range("E17").Value = '  -4.33%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '216.35'
$ws.Range("E18").Value = '  -5.51%  '
$ws.Range("D19").Value = '0.0₃0690'
$ws.Range("E19").Value = '  -4.11%  '
$ws.Range("E20").Value = '  -4.94%  '
$ws.Range("E21").Value = '  +0.38%  '
$ws.Range("E22").Value = '  -4.89%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.42'
$ws.Range("E23").Value = '  -6.28%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.98'
$ws.Range("E24").Value = '  -4.34%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.51'
$ws.Range("E25").Value = '  -1.06%  '
$ws.Range("E26").Value = '  +0.32%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.66'
$ws.Range("E27").Value = '  -3.25%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '14.99'
$ws.Range("E28").Value = '  -3.47%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.105'
$ws.Range("E29").Value = '  -4.89%  '
$ws.Range("E30").Value = '  -2.40%  '
$ws.Range("E31").Value = '  -3.94%  '
$ws.Range("E32").Value = '  -5.67%  '
$ws.Range("D33").Value = '1.361.64'
$ws.Range("E33").Value = '  -2.21%  '
$ws.Range("E34").Value = '  -6.18%  '
$ws.Range("E35").Value = '  -5.75%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.965'
$ws.Range("E36").Value = '  -5.08%  '
$ws.Range("E37").Value = '  -1.42%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0163'
$ws.Range("E38").Value = '  -4.29%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.535'
$ws.Range("E39").Value = '  -4.16%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.811'
$ws.Range("E40").Value = '  -4.51%  '
$ws.Range("E41").Value = '  +0.35%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.972'
$ws.Range("E42").Value = '  -4.20%  '
$ws.Range("B43").Value = 'MXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.17'
$ws.Range("E43").Value = '  +1.34%  '
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.76'
$ws.Range("E44").Value = '  -3.51%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '63.21'
$ws.Range("E45").Value = '  -3.85%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.26'
$ws.Range("E46").Value = '  -3.11%  '
$ws.Range("D47").Value = '1.710.85'
$ws.Range("E47").Value = '  -3.56%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.13'
$ws.Range("E48").Value = '  -1.82%  '
$ws.Range("D49").Value = '0.0₆0100'
$ws.Range("E49").Value = '  -3.22%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0965'
$ws.Range("E50").Value = '  -5.05%  '
$ws.Range("E51").Value = '  -1.65%  '
